# Apply updated dSF ("F" column) values per the commit:
# "repull data, push all data, mean calculation"
# The source data was repulled and the dSF column values were refreshed
# for each of the rows below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -2
    4  = -3
    5  = 1
    6  = 3
    7  = -2
    8  = -3
    9  = -1
    10 = -5
    11 = 3
    12 = 1
    14 = 2
    15 = 5
    17 = 1
    18 = 7
    20 = 1
    21 = 1
    22 = -3
    24 = -1
    25 = 3
    27 = 5
    28 = -2
    29 = -4
    30 = -2
    31 = -4
    32 = -3
    33 = -1
    34 = 1
    35 = 1
    36 = -2
    37 = 2
    39 = -1
    40 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
